# Trade #22 closed at 2026-02-17 13:18:27 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade (#22) for the MarketMaking strategy.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1199.24    # Current Capital
$summary.Cells.Item(4, 2).Value = -0.77      # Total P&L $
$summary.Cells.Item(5, 2).Value = -0.7       # Total P&L %
$summary.Cells.Item(6, 2).Value = 22         # Total Trades
$summary.Cells.Item(8, 2).Value = 14         # Losing Trades
$summary.Cells.Item(9, 2).Value = 31.82      # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(4, 3).Value = 99.23999999999999   # Capital
$status.Cells.Item(4, 4).Value = 22                  # Trades
$status.Cells.Item(4, 5).Value = -0.77               # P&L $
$status.Cells.Item(4, 6).Value = -0.76               # P&L %
$status.Cells.Item(4, 7).Value = 31.82               # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#22 / spreadsheet row 23) to a
# trade-log sheet ("All Trades" and the per-strategy "MarketMaking" sheet
# share the same layout).
# ---------------------------------------------------------------------
function Add-Trade23Row($sheet) {
    $sheet.Cells.Item(23, 1).Value = 22

    # Force the date column to stay literal text ("2026-02-17") instead of
    # being auto-parsed into a date serial number, matching the rest of the
    # column (which stores dates as plain text).
    $dateCell = $sheet.Cells.Item(23, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"

    $sheet.Cells.Item(23, 3).Value = "13:18:21"
    $sheet.Cells.Item(23, 4).Value = "MarketMaking"
    $sheet.Cells.Item(23, 5).Value = "DOWN"
    $sheet.Cells.Item(23, 6).Value = 0.66
    $sheet.Cells.Item(23, 7).Value = 0.636364
    $sheet.Cells.Item(23, 8).Value = "CLOSED"
    $sheet.Cells.Item(23, 9).Value = -3.5813
    $sheet.Cells.Item(23, 10).Value = -0.02
    $sheet.Cells.Item(23, 11).Value = 99.23999999999999
    $sheet.Cells.Item(23, 12).Value = 0
    $sheet.Cells.Item(23, 13).Value = 0
    $sheet.Cells.Item(23, 14).Value = 0.6
    $sheet.Cells.Item(23, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(23, 16).Value = "early_exit"
    $sheet.Cells.Item(23, 17).Value = 0.13
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade23Row $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade23Row $marketMaking
